$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: D-column cells whose new value is a numeric-looking string need to be
# forced to Text so Excel's "smart typing" doesn't silently convert them to a
# real number (and e.g. drop trailing zeros or switch to scientific notation).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.112.01"
$ws.Range("E2").Value = "  +0.93%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.127.46"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
Set-TextValue "D5" "599.95"
$ws.Range("E5").Value = "  -1.27%  "

# Row 6 - Solana
Set-TextValue "D6" "143.33"
$ws.Range("E6").Value = "  -0.62%  "

# Row 7
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.32%  "

# Row 8
$ws.Range("D8").Value = "3.133.21"
$ws.Range("E8").Value = "  +0.54%  "

# Row 9 (only E changes)
$ws.Range("E9").Value = "  +0.92%  "

# Row 10 (only E changes)
$ws.Range("E10").Value = "  +0.09%  "

# Row 11
Set-TextValue "D11" "5.41"
$ws.Range("E11").Value = "  +3.68%  "

# Row 12
Set-TextValue "D12" "0.470"
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
Set-TextValue "D13" "0.0000253"
$ws.Range("E13").Value = "  +1.51%  "

# Row 14
Set-TextValue "D14" "35.28"
$ws.Range("E14").Value = "  +0.35%  "

# Row 15
$ws.Range("D15").Value = "3.640.18"
$ws.Range("E15").Value = "  +0.24%  "

# Row 16 (only E changes)
$ws.Range("E16").Value = "  +2.97%  "

# Row 17
$ws.Range("D17").Value = "64.093.34"
$ws.Range("E17").Value = "  +0.75%  "

# Row 18
$ws.Range("D18").Value = "3.140.12"
$ws.Range("E18").Value = "  +0.86%  "

# Row 19
Set-TextValue "D19" "6.86"
$ws.Range("E19").Value = "  +1.13%  "

# Row 20
Set-TextValue "D20" "483.06"
$ws.Range("E20").Value = "  +1.87%  "

# Row 21
Set-TextValue "D21" "14.58"
$ws.Range("E21").Value = "  +0.25%  "

# Row 22
Set-TextValue "D22" "0.709"
$ws.Range("E22").Value = "  +0.96%  "

# Row 23
Set-TextValue "D23" "7.66"
$ws.Range("E23").Value = "  -0.79%  "

# Row 24
Set-TextValue "D24" "87.93"
$ws.Range("E24").Value = "  +5.57%  "

# Row 25
Set-TextValue "D25" "13.41"
$ws.Range("E25").Value = "  -0.34%  "

# Row 26
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.24%  "

# Row 27
Set-TextValue "D27" "2.73"
$ws.Range("E27").Value = "  -1.68%  "

# Row 28
Set-TextValue "D28" "8.34"
$ws.Range("E28").Value = "  -0.27%  "

# Row 29
Set-TextValue "D29" "7.12"
$ws.Range("E29").Value = "  +4.28%  "

# Row 30 (only E changes)
$ws.Range("E30").Value = "  -0.44%  "

# Row 31 (only E changes)
$ws.Range("E31").Value = "  -5.93%  "

# Row 32 (only E changes)
$ws.Range("E32").Value = "  +0.09%  "

# Row 33
Set-TextValue "D33" "26.90"
$ws.Range("E33").Value = "  +2.90%  "

# Row 34
Set-TextValue "D34" "2.66"
$ws.Range("E34").Value = "  -0.21%  "

# Row 35 (only E changes)
$ws.Range("E35").Value = "  -1.50%  "

# Row 36
Set-TextValue "D36" "6.02"
$ws.Range("E36").Value = "  +1.77%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0754"
$ws.Range("E37").Value = "  +0.28%  "

# Row 38
Set-TextValue "D38" "52.81"
$ws.Range("E38").Value = "  +0.25%  "

# Row 39
Set-TextValue "D39" "3.00"
$ws.Range("E39").Value = "  +2.43%  "

# Row 40
Set-TextValue "D40" "436.90"
$ws.Range("E40").Value = "  -3.43%  "

# Row 41 (only E changes)
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
Set-TextValue "D42" "0.119"
$ws.Range("E42").Value = "  +1.03%  "

# Row 43
Set-TextValue "D43" "8.24"
$ws.Range("E43").Value = "  -0.76%  "

# Row 44
$ws.Range("D44").Value = "2.860.66"
$ws.Range("E44").Value = "  +0.85%  "

# Row 45 - was TheGraph, now Fetch.AI
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D45" "2.23"
$ws.Range("E45").Value = "  -1.99%  "

# Row 46 - was Fetch.AI, now TheGraph
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D46" "0.259"
$ws.Range("E46").Value = "  -1.43%  "

# Row 47
Set-TextValue "D47" "2.43"
$ws.Range("E47").Value = "  +1.24%  "

# Row 48 (only E changes)
$ws.Range("E48").Value = "  -0.02%  "

# Row 49
Set-TextValue "D49" "25.83"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50 (only E changes)
$ws.Range("E50").Value = "  +0.83%  "

# Row 51
Set-TextValue "D51" "121.24"
$ws.Range("E51").Value = "  +2.22%  "
